# Sprint 2 documentation update
# - Sheet "Process": replace the placeholder "benefit/meaning" note (C4) with the
#   full project background / benefit text, which grows the row height.
# - Sheet "CSDL": just a cursor/selection move (handled at the end).
# - Sheet "họp regular Scrum": replace the "what I did" column (D) for the
#   sprint 2/3/4 regular-scrum rows with the actual sprint 2 task notes.
# - Sheet "workflow": just a cursor/selection move (handled at the end).

$wb = $excel.ActiveWorkbook

$wsProcess = $wb.Worksheets.Item("Process")
$wsCSDL = $wb.Worksheets.Item("CSDL")
$wsScrum = $wb.Worksheets.Item("h\u1ecdp regular Scrum")
$wsWorkflow = $wb.Worksheets.Item("workflow")

# ---------------------------------------------------------------------------
# Sheet "Process"
# ---------------------------------------------------------------------------

$benefitText = "T\u00ccNH TR\u1ea0NG CHUNG HI\u1ec6N NAY C\u1ee6A NHU C\u1ea6U T\u00ccM KI\u1ebeM VI\u1ec6C L\u00c0M `n" +
  "-Hi\u1ec7n t\u1ea1i vi\u1ec7c t\u00ecm ki\u1ebfm c\u00f4ng vi\u1ec7c c\u1ee7a sinh vi\u00ean \u0111ang l\u00e0 m\u1ed9t nhu c\u1ea7u thi\u1ebft y\u1ebfu.                             - R\u1ea5t nhi\u1ec1u nh\u00e0 tuy\u1ec3n d\u1ee5ng \u0111ang c\u00f3 nhu c\u1ea7u tuy\u1ec3n c\u00e1c sinh vi\u00ean v\u00e0o \u0111\u1ec3 \u0111\u00e0o t\u1ea1o v\u00e0 l\u00e0m vi\u1ec7c `n" +
  "L\u1ee2I \u00cdCH V\u00c0 \u00dd NGH\u0128A C\u1ee6A D\u1ef0 \u00c1N `n" +
  "- Gi\u00fap sinh vi\u00ean c\u00f3 \u0111\u01b0\u1ee3c vi\u1ec7c l\u00e0m ph\u00f9 h\u1ee3p       `n" +
  "-  Sinh vi\u00ean t\u00edch l\u0169y \u0111\u01b0\u1ee3c kinh nghi\u1ec7m khi \u0111i l\u00e0m th\u1ef1c t\u1ebf                                                               - Nh\u00e0 tuy\u1ec3n d\u1ee5ng c\u00f3 \u0111\u01b0\u1ee3c nh\u01b0ng nh\u00e2n vi\u00ean \u0111\u1ea7y nhi\u1ec7t huy\u1ebft v\u00e0 th\u00f4ng minh`n" +
  "V\u1eady n\u00f3 \u0111\u00e3 gi\u1ea3i quy\u1ebft \u0111\u01b0\u1ee3c nh\u1eefng v\u1ea5n \u0111\u1ec1 sau :`n" +
  "-  Sinh vi\u00ean hi\u1ec7n nay ch\u01b0a c\u00f3 vi\u1ec7c l\u00e0m ph\u00f9 h\u01a1p ( tr\u00e1i ng\u00e0nh ) ch\u01b0a c\u00f3 nhi\u1ec1u kinh nghi\u1ec7m, ch\u01b0a c\u00f3 m\u00f4i tr\u01b0\u1eddng r\u00e8n luy\u1ec7n ph\u00f9 h\u1ee3p`n" +
  "- Trong qu\u00e1 tr\u00ecnh h\u1ecdc t\u1eadp tr\u00ean tr\u01b0\u1eddng sinh vi\u00ean ch\u01b0a c\u00f3 kinh nghi\u1ec7m nhi\u1ec1u trong c\u00f4ng vi\u1ec7c, n\u00ean sau khi h\u1ecdc xong r\u1ea5t kh\u00f3 ki\u1ebfm \u0111\u01b0\u1ee3c vi\u1ec7c`n"

$wsProcess.Range("C4").Value = $benefitText
$wsProcess.Rows.Item(4).RowHeight = 85.5

# ---------------------------------------------------------------------------
# Sheet "h\u1ecdp regular Scrum" - sprint 2 "what I did" notes
# ---------------------------------------------------------------------------

$wsScrum.Range("D7").Value  = "Thi\u1ebft k\u1ebf database, vi\u1ebft user story cho ch\u1ee9c n\u0103ng \u0111\u0103ng tin"
$wsScrum.Range("D8").Value  = "Thi\u1ebft k\u1ebf giao di\u1ec7n website"
$wsScrum.Range("D9").Value  = "Vi\u1ebft userstory ch\u1ee9c n\u0103ng s\u1eeda , x\u00f3a tin"
$wsScrum.Range("D10").Value = "Vi\u1ebft userstory ch\u1ee9c n\u0103ng hi\u1ec3n th\u1ecb tin"
$wsScrum.Range("D11").Value = "60% ch\u1ee9c n\u0103ng \u0111\u0103ng tin"
$wsScrum.Range("D12").Value = "Ho\u00e0n th\u00e0nh giao di\u1ec7n"
$wsScrum.Range("D13").Value = "60% ch\u1ee9c n\u0103ng s\u1eeda tin x\u00f3a tin"
$wsScrum.Range("D14").Value = "60% ch\u1ee9c n\u0103ng hi\u1ec7n tin"
$wsScrum.Range("D15").Value = "t\u00edch h\u1ee3p v\u00e0o giao di\u1ec7n trang web"
$wsScrum.Range("D16").Value = "bug l\u1ed7i giao di\u1ec7n"
$wsScrum.Range("D17").Value = "t\u00edch h\u1ee3p v\u00e0o giao di\u1ec7n trang web"
$wsScrum.Range("D18").Value = "t\u00edch h\u1ee3p v\u00e0o giao di\u1ec7n trang web"

# ---------------------------------------------------------------------------
# Restore the cursor/selection state recorded for each sheet after the edit
# ---------------------------------------------------------------------------

$wsProcess.Activate()
$wsProcess.Application.ActiveWindow.ScrollRow = 19
$wsProcess.Range("C4").Select() | Out-Null

$wsCSDL.Activate()
$wsCSDL.Range("D6").Select() | Out-Null

$wsScrum.Activate()
$wsScrum.Range("F12").Select() | Out-Null

$wsWorkflow.Activate()
$wsWorkflow.Range("F23").Select() | Out-Null

$wsWorkflow.Activate()
